$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 13.115

$ws.Range("B4").Value = 6.37
$ws.Range("D4").Value = -8.099
$ws.Range("E4").Value = 12.425

$ws.Range("D5").Value = -8.581999999999999

$ws.Range("B6").Value = 7.295999999999999

$ws.Range("B7").Value = 7.013

$ws.Range("D8").Value = -8.204000000000001

$ws.Range("E9").Value = 13.022

$ws.Range("E11").Value = 12.852

$ws.Range("E14").Value = 13.06

$ws.Range("B16").Value = 6.575
$ws.Range("D16").Value = -8.296000000000001

$ws.Range("E18").Value = 12.596

$ws.Range("B20").Value = 6.09

$ws.Range("D22").Value = -8.191999999999998

$ws.Range("E25").Value = 12.791
